$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated "K" column (G) values - computed from strikeouts (Strike#) -> K
$kValues = @{
    2 = 0
    3 = 0
    4 = 2
    5 = 2
    6 = 2
    7 = 0
    8 = 1
    9 = 0
    10 = 1
    11 = 0
    12 = 3
    13 = 2
    14 = 0
    15 = 0
    16 = 0
    17 = 1
    18 = 2
    19 = 5
    20 = 0
    21 = 3
    22 = 0
    23 = 1
    25 = 1
    26 = 1
    27 = 2
    28 = 0
    29 = 0
    30 = 1
    31 = 2
    32 = 2
    33 = 1
    34 = 1
    35 = 1
    36 = 2
    37 = 1
    38 = 3
    39 = 4
    40 = 1
    41 = 0
    42 = 1
    43 = 1
    44 = 2
    45 = 2
    46 = 1
    47 = 1
    48 = 1
    49 = 1
    50 = 3
    51 = 5
    52 = 0
    53 = 1
    54 = 0
    55 = 0
    56 = 2
    57 = 2
    58 = 0
    59 = 2
    60 = 2
    61 = 1
    62 = 3
    63 = 1
    64 = 2
    65 = 3
    66 = 3
    67 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
